$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').NumberFormat = '@'
$ws.Range('D2').Value = '56.583.54'
$ws.Range('D2').Style = 'Normal'
$ws.Range('E2').Value = '  +5.33%  '
$ws.Range('D3').NumberFormat = '@'
$ws.Range('D3').Value = '2.510.14'
$ws.Range('D3').Style = 'Normal'
$ws.Range('E3').Value = '  +3.73%  '
$ws.Range('E4').Value = '  -0.16%  '
$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '490.91'
$ws.Range('D5').Style = 'Normal'
$ws.Range('E5').Value = '  +6.30%  '
$ws.Range('D6').NumberFormat = '@'
$ws.Range('D6').Value = '147.92'
$ws.Range('D6').Style = 'Normal'
$ws.Range('E6').Value = '  +12.76%  '
$ws.Range('E7').Value = '  -0.21%  '
$ws.Range('E8').Value = '  +6.52%  '
$ws.Range('D9').NumberFormat = '@'
$ws.Range('D9').Value = '2.532.21'
$ws.Range('D9').Style = 'Normal'
$ws.Range('E9').Value = '  +3.18%  '
$ws.Range('D10').NumberFormat = '@'
$ws.Range('D10').Value = '5.78'
$ws.Range('D10').Style = 'Normal'
$ws.Range('E10').Value = '  +8.69%  '
$ws.Range('D11').NumberFormat = '@'
$ws.Range('D11').Value = '0.0982'
$ws.Range('D11').Style = 'Normal'
$ws.Range('E11').Value = '  +3.64%  '
$ws.Range('E12').Value = '  +5.94%  '
$ws.Range('D13').NumberFormat = '@'
$ws.Range('D13').Value = '0.124'
$ws.Range('D13').Style = 'Normal'
$ws.Range('E13').Value = '  +1.64%  '
$ws.Range('D14').NumberFormat = '@'
$ws.Range('D14').Value = '2.942.18'
$ws.Range('D14').Style = 'Normal'
$ws.Range('E14').Value = '  +2.75%  '
$ws.Range('D15').NumberFormat = '@'
$ws.Range('D15').Value = '56.576.71'
$ws.Range('D15').Style = 'Normal'
$ws.Range('E15').Value = '  +5.34%  '
$ws.Range('D16').NumberFormat = '@'
$ws.Range('D16').Value = '21.33'
$ws.Range('D16').Style = 'Normal'
$ws.Range('E16').Value = '  +8.32%  '
$ws.Range('D17').NumberFormat = '@'
$ws.Range('D17').Value = '0.0000138'
$ws.Range('D17').Style = 'Normal'
$ws.Range('E17').Value = '  +3.28%  '
$ws.Range('D18').NumberFormat = '@'
$ws.Range('D18').Value = '2.529.50'
$ws.Range('D18').Style = 'Normal'
$ws.Range('E18').Value = '  +2.93%  '
$ws.Range('D19').NumberFormat = '@'
$ws.Range('D19').Value = '4.53'
$ws.Range('D19').Style = 'Normal'
$ws.Range('E19').Value = '  +8.39%  '
$ws.Range('E20').Value = '  +9.35%  '
$ws.Range('D21').NumberFormat = '@'
$ws.Range('D21').Value = '323.27'
$ws.Range('D21').Style = 'Normal'
$ws.Range('E21').Value = '  +4.12%  '
$ws.Range('D22').NumberFormat = '@'
$ws.Range('D22').Value = '1.00'
$ws.Range('D22').Style = 'Normal'
$ws.Range('E22').Value = '  +0.82%  '
$ws.Range('D23').NumberFormat = '@'
$ws.Range('D23').Value = '5.85'
$ws.Range('D23').Style = 'Normal'
$ws.Range('E23').Value = '  +9.49%  '
$ws.Range('D24').NumberFormat = '@'
$ws.Range('D24').Value = '58.88'
$ws.Range('D24').Style = 'Normal'
$ws.Range('E24').Value = '  +4.63%  '
$ws.Range('E25').Value = '  +8.31%  '
$ws.Range('E26').Value = '  +10.12%  '
$ws.Range('D27').NumberFormat = '@'
$ws.Range('D27').Value = '0.996'
$ws.Range('D27').Style = 'Normal'
$ws.Range('E27').Value = '  -1.49%  '
$ws.Range('D28').NumberFormat = '@'
$ws.Range('D28').Value = '2.619.49'
$ws.Range('D28').Style = 'Normal'
$ws.Range('E28').Value = '  +2.36%  '
$ws.Range('D29').NumberFormat = '@'
$ws.Range('D29').Value = '7.64'
$ws.Range('D29').Style = 'Normal'
$ws.Range('E29').Value = '  +6.38%  '
$ws.Range('D30').NumberFormat = '@'
$ws.Range('D30').Value = '0.0₃0805'
$ws.Range('D30').Style = 'Normal'
$ws.Range('E30').Value = '  +11.76%  '
$ws.Range('D31').NumberFormat = '@'
$ws.Range('D31').Value = '1.00'
$ws.Range('D31').Style = 'Normal'
$ws.Range('E31').Value = '  +0.05%  '
$ws.Range('D32').NumberFormat = '@'
$ws.Range('D32').Value = '149.47'
$ws.Range('D32').Style = 'Normal'
$ws.Range('E32').Value = '  -1.50%  '
$ws.Range('D33').NumberFormat = '@'
$ws.Range('D33').Value = '18.38'
$ws.Range('D33').Style = 'Normal'
$ws.Range('E33').Value = '  +4.40%  '
$ws.Range('E34').Value = '  +6.36%  '
$ws.Range('D35').NumberFormat = '@'
$ws.Range('D35').Value = '5.24'
$ws.Range('D35').Style = 'Normal'
$ws.Range('E35').Value = '  +4.72%  '
$ws.Range('D36').NumberFormat = '@'
$ws.Range('D36').Value = '1.15'
$ws.Range('D36').Style = 'Normal'
$ws.Range('E36').Value = '  +9.70%  '
$ws.Range('D37').NumberFormat = '@'
$ws.Range('D37').Value = '3.78'
$ws.Range('D37').Style = 'Normal'
$ws.Range('E37').Value = '  +7.10%  '
$ws.Range('E38').Value = '  +11.55%  '
$ws.Range('D39').NumberFormat = '@'
$ws.Range('D39').Value = '34.37'
$ws.Range('D39').Style = 'Normal'
$ws.Range('E39').Value = '  +2.78%  '
$ws.Range('D40').NumberFormat = '@'
$ws.Range('D40').Value = '3.57'
$ws.Range('D40').Style = 'Normal'
$ws.Range('E40').Value = '  +9.49%  '
$ws.Range('E41').Value = '  +4.11%  '
$ws.Range('D42').NumberFormat = '@'
$ws.Range('D42').Value = '0.0560'
$ws.Range('D42').Style = 'Normal'
$ws.Range('E42').Value = '  +6.26%  '
$ws.Range('D43').NumberFormat = '@'
$ws.Range('D43').Value = '0.993'
$ws.Range('D43').Style = 'Normal'
$ws.Range('E43').Value = '  -0.45%  '
$ws.Range('E44').Value = '  +8.75%  '
$ws.Range('E45').Value = '  +13.90%  '
$ws.Range('D46').NumberFormat = '@'
$ws.Range('D46').Value = '262.97'
$ws.Range('D46').Style = 'Normal'
$ws.Range('E46').Value = '  +19.49%  '
$ws.Range('D47').NumberFormat = '@'
$ws.Range('D47').Value = '0.0230'
$ws.Range('D47').Style = 'Normal'
$ws.Range('E47').Value = '  +5.45%  '
$ws.Range('D48').NumberFormat = '@'
$ws.Range('D48').Value = '0.0916'
$ws.Range('D48').Style = 'Normal'
$ws.Range('E48').Value = '  +5.56%  '
$ws.Range('E49').Value = '  -0.09%  '
$ws.Range('D50').NumberFormat = '@'
$ws.Range('D50').Value = '1.926.80'
$ws.Range('D50').Style = 'Normal'
$ws.Range('E50').Value = '  -1.35%  '
$ws.Range('D51').NumberFormat = '@'
$ws.Range('D51').Value = '17.75'
$ws.Range('D51').Style = 'Normal'
$ws.Range('E51').Value = '  +7.08%  '
